$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "awesome" -> "awesome (bonus)" (A12, same row still pairs with the
# "In the awesome branch..." refinement description in B12)
$ws.Range("A12").Value = "awesome (bonus)"

# Total score cell: 110 -> 100
$ws.Range("C16").Value = 100

# Scroll the window so row 4 is the top visible row, then select C16 -
# matches the sheetView's topLeftCell="A4" / selection activeCell="C16".
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C16").Select()
